# Rebuild the data table (rows 17-28) so that each worker's rows are grouped
# together (all YOVANIS ANTONIO VIGA OSORIO rows first, then all OSCAR JAVIER
# RONCANCIO VALBUENA rows), instead of the previous alternating order.
# The underlying records (person + period + amounts) are unchanged - only
# their row order/grouping in the sheet changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 17; Doc = "78756892";    Name = "YOVANIS ANTONIO VIGA OSORIO";        Period = "1708"; Mora = 29509;  Salario = 737717 },
    @{ Row = 18; Doc = "78756892";    Name = "YOVANIS ANTONIO VIGA OSORIO";        Period = "1707"; Mora = 29509;  Salario = 737717 },
    @{ Row = 19; Doc = "78756892";    Name = "YOVANIS ANTONIO VIGA OSORIO";        Period = "1706"; Mora = 29509;  Salario = 737717 },
    @{ Row = 20; Doc = "78756892";    Name = "YOVANIS ANTONIO VIGA OSORIO";        Period = "1705"; Mora = 29509;  Salario = 737717 },
    @{ Row = 21; Doc = "78756892";    Name = "YOVANIS ANTONIO VIGA OSORIO";        Period = "1704"; Mora = 29509;  Salario = 737717 },
    @{ Row = 22; Doc = "78756892";    Name = "YOVANIS ANTONIO VIGA OSORIO";        Period = "1702"; Mora = 29509;  Salario = 737717 },
    @{ Row = 23; Doc = "1032415619";  Name = "OSCAR JAVIER RONCANCIO VALBUENA";    Period = "1708"; Mora = 140000; Salario = 3500000 },
    @{ Row = 24; Doc = "1032415619";  Name = "OSCAR JAVIER RONCANCIO VALBUENA";    Period = "1707"; Mora = 140000; Salario = 3500000 },
    @{ Row = 25; Doc = "1032415619";  Name = "OSCAR JAVIER RONCANCIO VALBUENA";    Period = "1706"; Mora = 140000; Salario = 3500000 },
    @{ Row = 26; Doc = "1032415619";  Name = "OSCAR JAVIER RONCANCIO VALBUENA";    Period = "1705"; Mora = 140000; Salario = 3500000 },
    @{ Row = 27; Doc = "1032415619";  Name = "OSCAR JAVIER RONCANCIO VALBUENA";    Period = "1704"; Mora = 140000; Salario = 3500000 },
    @{ Row = 28; Doc = "1032415619";  Name = "OSCAR JAVIER RONCANCIO VALBUENA";    Period = "1702"; Mora = 140000; Salario = 3500000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc      # C - N Doc Trabajador
    $ws.Cells.Item($r.Row, 4).Value = $r.Name     # D - Nombre Trabajador
    $ws.Cells.Item($r.Row, 5).Value = $r.Period   # E - Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora     # F - Valor Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario  # G - Salario Basico
}
